$wb = $excel.ActiveWorkbook

# --- Update existing sheets with corrected figures ---
$ws1 = $wb.Worksheets.Item("Stock log")

# Row 2 (Atlantica Sustainable Infrastru, AY)
$ws1.Range("F2").Value = 21.69
$ws1.Range("J2").Value = 14901.03
$ws1.Range("K2").Value = -7.9
$ws1.Range("M2").Value = -4.28

# Row 3 (Atlantica Sustainable Infrastru, AY)
$ws1.Range("F3").Value = 21.69
$ws1.Range("J3").Value = 8676
$ws1.Range("K3").Value = -12.04
$ws1.Range("M3").Value = -8.75

# Row 4 (Schwab US Dividend Equity ETF, SCHD)
$ws1.Range("F4").Value = 76.28
$ws1.Range("J4").Value = 74983.24
$ws1.Range("K4").Value = 0.83
$ws1.Range("M4").Value = 2.73

# Row 5 (Atlantica Sustainable Infrastru, AY - Sell)
$ws1.Range("F5").Value = 21.69
$ws1.Range("J5").Value = -2169
$ws1.Range("K5").Value = -5.7
$ws1.Range("M5").Value = -3.86

$ws2 = $wb.Worksheets.Item("Portfolio Summary")

# Row 2 (AY)
$ws2.Range("D2").Value = 21408.03
$ws2.Range("G2").Value = -9.85
$ws2.Range("H2").Value = -6.19

# Row 3 (SCHD)
$ws2.Range("D3").Value = 74983.24
$ws2.Range("G3").Value = 0.83
$ws2.Range("H3").Value = 2.73

# --- Add new "Total Return" sheet comparing portfolio performance with indexes ---
$count = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($count)
$ws3 = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$ws3.Name = "Total Return"

# Copy header formatting (bold, bordered, centered) from the "Stock log" sheet
$ws1.Range("A1:D1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats
$ws1.Range("A1").Copy()
$ws3.Range("A2:A4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws3.Range("B1").Value = "Portfolio"
$ws3.Range("C1").Value = "Price Return, %"
$ws3.Range("D1").Value = "Total Return, %"

$ws3.Range("A2").Value = 0
$ws3.Range("B2").Value = "My Portfolio"
$ws3.Range("C2").Value = -1.75
$ws3.Range("D2").Value = 0.61

$ws3.Range("A3").Value = 1
$ws3.Range("B3").Value = "Global X Super Dividend ETF"
$ws3.Range("C3").Value = 10
$ws3.Range("D3").Value = -0.88

$ws3.Range("A4").Value = 2
$ws3.Range("B4").Value = "S&P 500"
$ws3.Range("C4").Value = 10
$ws3.Range("D4").Value = 26.51

$ws1.Select()
